$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns at D:E, shifting existing D:K -> F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from the (now-shifted) old D:E column (now at F:G) into new D:E
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Write updated financial data for row ranges D:M (rows 7-35, 38-77, 80-102) ---
function Set-RowValues($rowNum, $values) {
    $arr = New-Object "object[,]" 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $values[$i] }
    $ws.Range("D$rowNum`:M$rowNum").Value2 = $arr
}

Set-RowValues 7 @(43463, 43379, 43295, 43211, 43099, 43015, 42931, 42847, 42735, 42651)
Set-RowValues 8 @(1896800, 1886700, 4281000, 2385100, 1885500, 1868400, 4209900, 2353700, 1828200, 1800100)
Set-RowValues 9 @(1651400, 1630600, 3672200, 2041900, 1630700, 1606700, 3581500, 1996300, 1568900, 1544800)
Set-RowValues 10 @(245400, 256100, 608800, 343200, 254800, 261700, 628400, 357400, 259300, 255300)
Set-RowValues 11 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 12 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues 13 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 14 @(33700, 800, 8000, 8400, 4300, 227000, 5600, 5000, 11400, 5100)
Set-RowValues 15 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 16 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 17 @(1908700, 1859900, 4225500, 2359300, 1867100, 2062300, 4141500, 2324100, 1803800, 1770200)
Set-RowValues 18 @(-11900, 26800, 55500, 25800, 18400, -193900, 68400, 29600, 24400, 29900)
Set-RowValues 19 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 20 @(0, 200, 500, 200, 300, 100, 300, 200, 100, 100)
Set-RowValues 21 @(7900, 46600, 100900, 51500, 36700, -173800, 115100, 55600, 43200, 48400)
Set-RowValues 22 @(7700, 7100, 15700, 8800, 6200, 6100, 13000, 7300, 4400, 4400)
Set-RowValues 23 @(-19500, 19900, 40300, 17200, 12500, -199900, 55800, 22400, 20100, 25600)
Set-RowValues 24 @(-5500, 2900, 10000, 4800, 3800, -76400, 19600, 7400, 7300, 8900)
Set-RowValues 25 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 26 @(-14000, 17100, 30300, 12400, 8700, -123500, 36100, 15100, 12800, 16700)
Set-RowValues 27 @(-13700, 16600, 29600, 12200, 8100, -121400, 35600, 14800, 12600, 16400)
Set-RowValues 28 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 29 @(0, 400, -200, -100, 25900, -100, -100, 0, 0, -100)
Set-RowValues 30 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 31 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 32 @(0, -200, -500, -200, -300, -100, -300, -200, -100, -100)
Set-RowValues 33 @(-13700, 17100, 29500, 12100, 34000, -121400, 35500, 14700, 12600, 16400)
Set-RowValues 34 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 35 @(-13700, 17100, 29500, 12100, 34000, -121400, 35500, 14700, 12600, 16400)
Set-RowValues 38 @(43463, 43379, 43295, 43211, 43099, 43015, 42931, 42847, 42735, 42651)
Set-RowValues 39 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 40 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 41 @(18600, 20700, 15900, 17200, 15700, 13200, 22700, 19500, 24400, 26400)
Set-RowValues 42 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 43 @(346300, 364000, 355100, 330300, 344100, 370500, 349300, 342400, 291600, 322000)
Set-RowValues 44 @(553800, 592200, 562400, 577500, 597200, 598500, 555600, 539900, 539900, 561800)
Set-RowValues 45 @(82500, 52000, 52400, 60200, 47400, 33400, 32900, 42900, 37700, 29600)
Set-RowValues 46 @(1001100, 1028800, 985800, 985200, 1004300, 1015600, 960500, 944700, 893500, 939700)
Set-RowValues 47 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 48 @(579100, 577300, 581800, 582400, 600200, 588400, 621600, 628000, 559700, 570700)
Set-RowValues 49 @(307600, 308900, 309800, 311200, 313100, 314000, 496700, 498900, 322700, 322700)
Set-RowValues 50 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 51 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 52 @(84200, 139100, 133400, 138100, 138200, 115800, 119800, 109000, 154400, 160700)
Set-RowValues 53 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 54 @(1971900, 2054000, 2010800, 2016900, 2055800, 2033800, 2198500, 2180600, 1930300, 1993900)
Set-RowValues 55 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 56 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 57 @(357800, 411400, 364700, 356500, 377000, 440600, 394300, 370700, 372400, 398900)
Set-RowValues 58 @(18300, 8100, 7800, 7600, 9200, 19400, 19000, 17400, 17400, 19000)
Set-RowValues 59 @(100400, 98600, 105900, 100400, 108400, 100000, 101500, 101400, 116100, 107100)
Set-RowValues 60 @(476500, 518100, 478400, 464500, 494600, 560000, 514800, 489500, 506000, 525100)
Set-RowValues 61 @(679800, 694900, 702900, 733400, 740800, 651500, 641300, 658300, 413700, 475400)
Set-RowValues 62 @(99700, 104700, 105100, 107000, 98500, 116300, 193400, 191400, 185300, 177900)
Set-RowValues 63 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 64 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 65 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 66 @(1256000, 1317700, 1286300, 1304900, 1333800, 1327800, 1349500, 1339200, 1104900, 1178300)
Set-RowValues 67 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 68 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 69 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 70 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 71 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 72 @(247600, 268100, 257100, 245800, 240000, 208800, 338400, 323600, 314900, 307600)
Set-RowValues 73 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 74 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 75 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 76 @(715900, 736400, 724500, 712000, 722000, 706000, 849100, 841500, 825400, 815600)
Set-RowValues 77 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 80 @(43463, 43379, 43295, 43211, 43099, 43015, 42931, 42847, 42735, 42651)
Set-RowValues 81 @(-13700, 17100, 29500, 12100, 34000, -121400, 35500, 14700, 12600, 16400)
Set-RowValues 82 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 83 @(19700, 19600, 44900, 25600, 18000, 20000, 46400, 25900, 18700, 18400)
Set-RowValues 84 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 85 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 86 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 87 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 88 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 89 @(29100, 38200, 104200, 60300, -18800, 33100, 38400, -10200, 75700, 23800)
Set-RowValues 90 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 91 @(-18900, -18000, -34600, -21000, -15600, -17500, -37800, -19500, -32400, -31800)
Set-RowValues 92 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 93 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 94 @(-18600, -17400, -28100, -20900, -38200, -29400, -247800, -232900, -15700, -17000)
Set-RowValues 95 @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
Set-RowValues 96 @(-6500, -6500, -13000, -6500, -6100, -6100, -12500, -6300, -5600, -5600)
Set-RowValues 97 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 98 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 99 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 100 @(-12500, -16000, -75800, -37900, 59600, -13400, 207800, 238300, -62000, -4200)
Set-RowValues 101 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 102 @(-2100, 4800, 200, 1500, 2500, -9600, -1600, -4800, -2000, 2600)
